$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.816.91"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "2.295.37"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("D5").Value = "'299.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.47%  "
$ws.Range("D6").Value = "'97.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.68%  "
$ws.Range("E7").Value = "  +1.01%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.504"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.99%  "
$ws.Range("D10").Value = "'35.81"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.77%  "
$ws.Range("E11").Value = "  -0.58%  "
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").Value = "'17.67"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.67%  "
$ws.Range("D14").Value = "'6.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.78%  "
$ws.Range("D15").Value = "2.651.44"
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("D16").Value = "2.306.92"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").Value = "'0.774"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.65%  "
$ws.Range("D18").Value = "42.779.56"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("E19").Value = "  -5.58%  "
$ws.Range("D20").Value = "0.0₃0905"
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("E21").Value = "  -2.30%  "
$ws.Range("D22").Value = "'67.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.64%  "
$ws.Range("D23").Value = "'241.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("E24").Value = "  -1.45%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'2.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.11%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "'4.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("D28").Value = "'25.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.52%  "
$ws.Range("D29").Value = "'165.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.06%  "
$ws.Range("E30").Value = "  -1.41%  "
$ws.Range("E31").Value = "  -2.09%  "
$ws.Range("D32").Value = "'32.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.95%  "
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D34").Value = "'4.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.09%  "
$ws.Range("E35").Value = "  -3.93%  "
$ws.Range("D36").Value = "'17.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.65%  "
$ws.Range("E37").Value = "  -0.95%  "
$ws.Range("D38").Value = "'0.0685"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.87%  "
$ws.Range("E39").Value = "  -2.34%  "
$ws.Range("D40").Value = "'1.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.51%  "
$ws.Range("E41").Value = "  -1.85%  "
$ws.Range("D42").Value = "'0.109"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("D43").Value = "2.013.61"
$ws.Range("E43").Value = "  +0.97%  "
$ws.Range("E44").Value = "  -2.16%  "
$ws.Range("D45").Value = "'10.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.45%  "
$ws.Range("D46").Value = "'2.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.36%  "
$ws.Range("D47").Value = "'17.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.03%  "
$ws.Range("D48").Value = "'2.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.42%  "
$ws.Range("D49").Value = "2.513.67"
$ws.Range("E49").Value = "  -1.29%  "
$ws.Range("D50").Value = "'52.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.84%  "
$ws.Range("D51").Value = "'2.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.67%  "
